$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "dateTime, fmt"
$ws.Range("B5").Value = 46001.524259259262
$ws.Range("B5").NumberFormat = "dd/mm/yyyy\ hh:mm:ss"

$ws.Range("B5").Select()
